$wb = $excel.ActiveWorkbook

# ===== ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 184.07143
$ws.Range("I33").Value = 117.333336
$ws.Range("J33").Value = 304.2
$ws.Range("K33").Value = 117.333336
$ws.Range("L33").Value = 304.2
$ws.Range("M33").Value = 111.666664
$ws.Range("N33").Value = -762.2
$ws.Range("H62").Value = 3715.8333
$ws.Range("I62").Value = 3715.8333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3715.8333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3091.8333
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3715.8333
$ws.Range("I65").Value = 3715.8333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 18579.1665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -15459.1665
$ws.Range("N65").ClearContents()
$ws.Range("H100").Value = 1482.6522
$ws.Range("I100").Value = 1208.6364
$ws.Range("J100").Value = 1733.8334
$ws.Range("K100").Value = 1208.6364
$ws.Range("L100").Value = 1733.8334
$ws.Range("M100").Value = -667.6364000000001
$ws.Range("N100").Value = -2815.8334
$ws.Range("H105").Value = 81396
$ws.Range("J105").Value = 81396
$ws.Range("L105").Value = 81396
$ws.Range("N105").Value = -88384
$ws.Range("H107").Value = 788.6896400000001
$ws.Range("I107").Value = 840.3125
$ws.Range("J107").Value = 725.1539
$ws.Range("K107").Value = 840.3125
$ws.Range("L107").Value = 725.1539
$ws.Range("M107").Value = 1079.6875
$ws.Range("N107").Value = -4565.1539
$ws.Range("H113").Value = 2053.6365
$ws.Range("I113").Value = 1998
$ws.Range("J113").Value = 2100
$ws.Range("K113").Value = 1998
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1256
$ws.Range("N113").Value = -8608
$ws.Range("H116").Value = 2468.1667
$ws.Range("I116").Value = 2561.8
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2561.8
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 880.1999999999998
$ws.Range("N116").Value = -8884
$ws.Range("H132").Value = 2838.4
$ws.Range("I132").Value = 2378.7673
$ws.Range("J132").Value = 5661.857
$ws.Range("K132").Value = 7136.3019
$ws.Range("L132").Value = 16985.571
$ws.Range("M132").Value = -4606.3019
$ws.Range("N132").Value = -22045.571

# ===== ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2388.2666
$ws.Range("I2").Value = 2247.6365
$ws.Range("J2").Value = 2775
$ws.Range("K2").Value = 2247.6365
$ws.Range("L2").Value = 2775
$ws.Range("M2").Value = -2134.6365
$ws.Range("N2").Value = -3001
$ws.Range("H92").Value = 34708
$ws.Range("J92").Value = 34708
$ws.Range("L92").Value = 34708
$ws.Range("N92").Value = -39700
$ws.Range("H97").Value = 1193.091
$ws.Range("I97").Value = 725
$ws.Range("J97").Value = 2012.25
$ws.Range("K97").Value = 725
$ws.Range("L97").Value = 2012.25
$ws.Range("M97").Value = -229
$ws.Range("N97").Value = -3004.25
$ws.Range("H110").Value = 3631
$ws.Range("I110").Value = 2011
$ws.Range("J110").Value = 4171
$ws.Range("K110").Value = 2011
$ws.Range("L110").Value = 4171
$ws.Range("M110").Value = 34
$ws.Range("N110").Value = -8261
$ws.Range("H116").Value = 2388.2666
$ws.Range("I116").Value = 2247.6365
$ws.Range("J116").Value = 2775
$ws.Range("K116").Value = 2247.6365
$ws.Range("L116").Value = 2775
$ws.Range("M116").Value = 46.36349999999993
$ws.Range("N116").Value = -7363
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ===== BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2388.2666
$ws.Range("I3").Value = 2247.6365
$ws.Range("J3").Value = 2775
$ws.Range("K3").Value = 2247.6365
$ws.Range("L3").Value = 2775
$ws.Range("M3").Value = -2133.6365
$ws.Range("N3").Value = -3003
$ws.Range("H99").Value = 3582.4443
$ws.Range("I99").Value = 2506.6667
$ws.Range("J99").Value = 4120.3335
$ws.Range("K99").Value = 2506.6667
$ws.Range("L99").Value = 4120.3335
$ws.Range("M99").Value = -1008.6667
$ws.Range("N99").Value = -7116.3335
$ws.Range("H133").Value = 30760
$ws.Range("J133").Value = 30760
$ws.Range("L133").Value = 30760
$ws.Range("N133").Value = -40880

# ===== CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4572
$ws.Range("N41").ClearContents()
$ws.Range("H62").Value = 2590.111
$ws.Range("I62").Value = 1999.6666
$ws.Range("J62").Value = 2885.3333
$ws.Range("K62").Value = 1999.6666
$ws.Range("L62").Value = 2885.3333
$ws.Range("M62").Value = -1375.6666
$ws.Range("N62").Value = -4133.3333
$ws.Range("H65").Value = 2590.111
$ws.Range("I65").Value = 1999.6666
$ws.Range("J65").Value = 2885.3333
$ws.Range("K65").Value = 9998.333000000001
$ws.Range("L65").Value = 14426.6665
$ws.Range("M65").Value = -6878.333000000001
$ws.Range("N65").Value = -20666.6665
$ws.Range("H105").Value = 779.2143
$ws.Range("I105").Value = 677.6667
$ws.Range("J105").Value = 962
$ws.Range("K105").Value = 677.6667
$ws.Range("L105").Value = 962
$ws.Range("M105").Value = 1069.3333
$ws.Range("N105").Value = -4456

# ===== CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4042.516
$ws.Range("I109").Value = 643.6
$ws.Range("J109").Value = 4696.154
$ws.Range("K109").Value = 1930.8
$ws.Range("L109").Value = 14088.462
$ws.Range("M109").Value = -890.8000000000002
$ws.Range("N109").Value = -16168.462
$ws.Range("H137").Value = 3213.087
$ws.Range("J137").Value = 3650
$ws.Range("L137").Value = 10950
$ws.Range("N137").Value = -21150

# ===== GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2894.6453
$ws.Range("I80").Value = 2740.6667
$ws.Range("J80").Value = 3039
$ws.Range("K80").Value = 2740.6667
$ws.Range("L80").Value = 3039
$ws.Range("M80").Value = -1742.6667
$ws.Range("N80").Value = -5035
$ws.Range("H83").Value = 2894.6453
$ws.Range("I83").Value = 2740.6667
$ws.Range("J83").Value = 3039
$ws.Range("K83").Value = 13703.3335
$ws.Range("L83").Value = 15195
$ws.Range("M83").Value = -8711.333500000001
$ws.Range("N83").Value = -25179

# ===== LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3155.76
$ws.Range("I40").Value = 2433.6667
$ws.Range("J40").Value = 3822.3076
$ws.Range("K40").Value = 2433.6667
$ws.Range("L40").Value = 3822.3076
$ws.Range("M40").Value = -2297.6667
$ws.Range("N40").Value = -4094.3076
$ws.Range("H61").Value = 1690.5
$ws.Range("I61").Value = 1642.4546
$ws.Range("J61").Value = 1866.6666
$ws.Range("K61").Value = 1642.4546
$ws.Range("L61").Value = 1866.6666
$ws.Range("M61").Value = -1440.4546
$ws.Range("N61").Value = -2270.6666
$ws.Range("H82").Value = 2886.9473
$ws.Range("I82").Value = 2627.4546
$ws.Range("J82").Value = 3243.75
$ws.Range("K82").Value = 2627.4546
$ws.Range("L82").Value = 3243.75
$ws.Range("M82").Value = -2266.4546
$ws.Range("N82").Value = -3965.75
$ws.Range("H85").Value = 2886.9473
$ws.Range("I85").Value = 2627.4546
$ws.Range("J85").Value = 3243.75
$ws.Range("K85").Value = 2627.4546
$ws.Range("L85").Value = 3243.75
$ws.Range("M85").Value = -1379.4546
$ws.Range("N85").Value = -5739.75
$ws.Range("H113").Value = 1690.5
$ws.Range("I113").Value = 1642.4546
$ws.Range("J113").Value = 1866.6666
$ws.Range("K113").Value = 1642.4546
$ws.Range("L113").Value = 1866.6666
$ws.Range("M113").Value = 527.5454
$ws.Range("N113").Value = -6206.6666
$ws.Range("H122").Value = 3616.1072
$ws.Range("I122").Value = 3198.682
$ws.Range("J122").Value = 5146.6665
$ws.Range("K122").Value = 9596.045999999998
$ws.Range("L122").Value = 15439.9995
$ws.Range("M122").Value = -7146.045999999998
$ws.Range("N122").Value = -20339.9995
$ws.Range("H137").Value = 38000
$ws.Range("J137").Value = 38000
$ws.Range("L137").Value = 38000
$ws.Range("N137").Value = -48200
$ws.Range("H139").Value = 44766
$ws.Range("J139").Value = 44766
$ws.Range("L139").Value = 44766
$ws.Range("N139").Value = -55046

# ===== WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 728.8461
$ws.Range("I113").Value = 713.1667
$ws.Range("J113").Value = 742.2857
$ws.Range("K113").Value = 2139.5001
$ws.Range("L113").Value = 2226.8571
$ws.Range("M113").Value = 30.4998999999998
$ws.Range("N113").Value = -6566.8571
